# Update "想去人数" (F) and "最低票价" (G) figures that changed between
# scrapes for two events that appear in both the "展览" sheet and the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (events listing) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8114
$ws1.Range("F5").Value = 5915
$ws1.Range("F6").Value = 501
$ws1.Range("F10").Value = 296
$ws1.Range("F11").Value = 457
$ws1.Range("G11").Value = 64.8

# --- Sheet "全部类型" (all event types combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8114
$ws4.Range("F5").Value = 5915
$ws4.Range("F6").Value = 501
$ws4.Range("F10").Value = 296
$ws4.Range("F15").Value = 457
$ws4.Range("G15").Value = 64.8
